$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 = "Week 3" (column A header is per-row week).
# F4 (Husam Alghamdi's entry for Week 3) gets an appended sentence about the
# activity diagram and planning poker.
$ws.Range("F4").Value = "Downloaded softwares to work in project in accordance with team plan, which are: ReactJS,NodeJS, MongoDB and Visual Studio Code. Contributed to User-Stories and Created project page on taiga. Added Activity digram for Agora and contribute to playing planning poker with the rest of the team."

# Update the active selection to reflect the cell that was last edited.
$ws.Range("F4").Select()
